$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 1
$ws.Range("A1").Value = "x"
$ws.Range("B1").Value = "y"

# Update the view: scroll back to top-left A1, and change selection to E7
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E7").Select()
